$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "309.02"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.21%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "36.99"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.68%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.120"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.43%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07793"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.67%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.397"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.21%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.301"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.47%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.56%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.25%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9233"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.42%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1098"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-9.51%"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.37%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08809"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-4.21%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03285"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.90%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09575"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.62%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001380"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.02%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006204"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "4.80%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.391"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-4.07%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3451"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.19%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.385"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "20.62%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.52%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2372"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-8.38%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04343"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.92%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-3.78%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004276"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.52%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "7.98%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002903"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02154"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "2.08%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04949"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-3.29%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007590"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.14%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1354"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.63%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008499"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-6.77%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002072"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "1.14%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008607"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-0.05%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006585"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.50%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.14%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "13.80%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001444"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "20.45%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.14%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.14%"
